# S5_N7_PCR-run1.xlsx -- "updated data; trying multi-file input"
#
# The "SYBR" sheet's G column ("End RFU") held computed readout values.
# This edit:
#   * drops the G-column value entirely for wells in rows 3-30 (the cell
#     is removed from the row, not just blanked), and
#   * overwrites the G-column value with the sentinel 10000 for the wells
#     in rows 32-94 that previously carried a computed (non-10000) value.
# It also moves the frozen-pane scroll position / active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYBR")

# Rows 3-30: remove the G-column cell outright.
$clearRows = 3..30
foreach ($r in $clearRows) {
    $ws.Cells.Item($r, 7).Clear()
}

# Rows whose G-column value becomes 10000.
$setRows = @(32,33,34,35,36,37,41,42,43,44,45,46,47,48,49,56,57,58,59,60,61,65,66,67,68,69,70,71,72,73,77,78,79,80,81,82,83,84,85,89,90,91,92,93,94)
foreach ($r in $setRows) {
    $ws.Cells.Item($r, 7).Value = 10000
}

# Move the active cell / scroll position within the frozen pane.
$ws.Range("M82").Select() | Out-Null
